$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($rangeAddr, $val) {
    $c = $ws.Range($rangeAddr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextCell 'D2' '54.216.01'
$ws.Range('E2').Value = '  -8.04%  '
Set-TextCell 'D3' '2.869.94'
$ws.Range('E3').Value = '  -11.14%  '
$ws.Range('E4').Value = '  +0.03%  '
Set-TextCell 'D5' '474.17'
$ws.Range('E5').Value = '  -12.25%  '
Set-TextCell 'D6' '126.14'
$ws.Range('E6').Value = '  -7.49%  '
$ws.Range('E7').Value = '  -0.01%  '
Set-TextCell 'D8' '2.865.71'
$ws.Range('E8').Value = '  -11.30%  '
$ws.Range('E9').Value = '  -12.08%  '
Set-TextCell 'D10' '6.65'
$ws.Range('E10').Value = '  -12.93%  '
Set-TextCell 'D11' '0.0974'
$ws.Range('E11').Value = '  -15.20%  '
$ws.Range('E12').Value = '  -15.53%  '
$ws.Range('E13').Value = '  -3.97%  '
Set-TextCell 'D14' '3.362.61'
$ws.Range('E14').Value = '  -11.15%  '
Set-TextCell 'D15' '22.76'
$ws.Range('E15').Value = '  -12.46%  '
Set-TextCell 'D16' '54.164.75'
$ws.Range('E16').Value = '  -8.19%  '
Set-TextCell 'D17' '2.878.67'
$ws.Range('E17').Value = '  -10.81%  '
Set-TextCell 'D18' '0.0000135'
$ws.Range('E18').Value = '  -14.85%  '
Set-TextCell 'D19' '5.21'
$ws.Range('E19').Value = '  -11.92%  '
$ws.Range('E20').Value = '  -13.21%  '
$ws.Range('E21').Value = '  -13.60%  '
Set-TextCell 'D22' '308.77'
$ws.Range('E22').Value = '  -14.97%  '
Set-TextCell 'D23' '0.997'
$ws.Range('E23').Value = '  -0.34%  '
Set-TextCell 'D24' '0.448'
$ws.Range('E24').Value = '  -14.12%  '
Set-TextCell 'D25' '59.68'
$ws.Range('E25').Value = '  -15.51%  '
$ws.Range('E26').Value = '  -0.37%  '
Set-TextCell 'D27' '0.154'
$ws.Range('E27').Value = '  -10.58%  '
$ws.Range('E28').Value = '  -0.09%  '
Set-TextCell 'D29' '0.0₃0818'
$ws.Range('E29').Value = '  -15.93%  '
$ws.Range('E30').Value = '  -12.11%  '
Set-TextCell 'D31' '1.14'
$ws.Range('E31').Value = '  -6.29%  '
Set-TextCell 'D32' '6.19'
$ws.Range('E32').Value = '  -13.05%  '
Set-TextCell 'D33' '19.07'
$ws.Range('E33').Value = '  -13.04%  '
$ws.Range('E34').Value = '  -16.42%  '
$ws.Range('E35').Value = '  -13.99%  '
Set-TextCell 'D36' '138.27'
$ws.Range('E36').Value = '  -14.28%  '
Set-TextCell 'D37' '5.45'
$ws.Range('E37').Value = '  -15.31%  '
Set-TextCell 'D38' '1.21'
$ws.Range('E38').Value = '  -15.98%  '
Set-TextCell 'D39' '23.09'
$ws.Range('E39').Value = '  -12.41%  '
Set-TextCell 'D40' '0.0621'
$ws.Range('E40').Value = '  -12.52%  '
Set-TextCell 'D41' '2.895.80'
$ws.Range('E41').Value = '  -11.18%  '
Set-TextCell 'D42' '0.999'
$ws.Range('E42').Value = '  -0.13%  '
Set-TextCell 'D43' '35.34'
$ws.Range('E43').Value = '  -14.07%  '
Set-TextCell 'D44' '0.960'
$ws.Range('E44').Value = '  -13.70%  '
Set-TextCell 'D45' '0.598'
$ws.Range('E45').Value = '  -16.55%  '
Set-TextCell 'D46' '3.42'
$ws.Range('E46').Value = '  -15.36%  '
$ws.Range('E47').Value = '  -12.70%  '
Set-TextCell 'D48' '2.058.17'
$ws.Range('E48').Value = '  -10.64%  '
$ws.Range('E49').Value = '  -15.42%  '
Set-TextCell 'D50' '17.88'
$ws.Range('E50').Value = '  -14.60%  '
$ws.Range('E51').Value = '  -12.07%  '
